$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.187.78"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").Value = "1.814.38"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3748"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8668"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").Value = "1.825.85"
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.647"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.381"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07102"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008740"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").Value = "27.183.41"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.313"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("D24").Value = "2.047.88"
$ws.Range("E24").Value = "  -5.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.932"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.232"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.267"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08887"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7728"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.171"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.514"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.924"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  +3.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01959"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05240"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.244"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.45%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.917"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.05%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.366"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +20.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5277"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1681"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.572"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5026"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.002"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.667"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06322"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "
